$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46074
$ws.Range("B2").Value = 32.97
$ws.Range("C2").Value = 27.75
$ws.Range("D2").Value = 24.1
$ws.Range("E2").Value = 23
$ws.Range("F2").Value = 20.01
$ws.Range("G2").Value = 23
$ws.Range("H2").Value = 27
$ws.Range("I2").Value = 28.71
$ws.Range("J2").Value = 25.96
$ws.Range("K2").Value = 10.21
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 0.59
$ws.Range("N2").Value = 0.58
$ws.Range("O2").Value = 0.59
$ws.Range("P2").Value = 0.55
$ws.Range("Q2").Value = 0.57
$ws.Range("R2").Value = 0.9399999999999999
$ws.Range("S2").Value = 9.630000000000001
$ws.Range("T2").Value = 29.41
$ws.Range("U2").Value = 53.36
$ws.Range("V2").Value = 91.45
$ws.Range("W2").Value = 86.37
$ws.Range("X2").Value = 44.96
$ws.Range("Y2").Value = 32.41
$ws.Range("Z2").Value = 24.8
$ws.Range("AB2").Value = 63.8
$ws.Range("AD2").Value = 88.91
$ws.Range("AF2").Value = 41.38
$ws.Range("AG2").Value = "2h-17h"
